$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.418.98"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.851.99"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.02"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6290"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07622"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2919"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.68"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07760"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.033"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6811"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001051"
$ws.Range("E14").Value = "  -5.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.16"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.138"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "29.421.20"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.23"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.36"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.499"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "159.12"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1389"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.70"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.457"
$ws.Range("E27").Value = "  +10.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.475"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05610"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.115"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.076"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.838"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.161"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7030"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.592"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").Value = "1.237.20"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01806"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.736"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.428"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9048"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.59"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.64"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.217"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4006"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.005"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.687"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1152"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("E49").Value = "  -4.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05701"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4629"
$ws.Range("E51").Value = "  +0.05%  "
